# Team Attendance [B8-G1] 10-Aug.
# Add a new attendance row (row 8) for 10-Aug-2023 (serial 45148) with
# PRESENT for columns B:E and ABSENT for columns F:I, plus "No Response"
# review comments on the ABSENT cells, matching the existing sheet pattern.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Date cell, formatted the same way as the other date cells in the column (d-mmm)
$ws.Range("A8").Value = 45148
$ws.Range("A8").NumberFormat = "d-mmm"

# Attendance values
$ws.Range("B8:E8").Value = "PRESENT"
$ws.Range("F8:I8").Value = "ABSENT"

# Comments on the ABSENT cells
$ws.Range("F8").AddComment("RENUKA:`nNo Response") | Out-Null
$ws.Range("G8").AddComment("RENUKA:`nNo Response") | Out-Null
$ws.Range("H8").AddComment("RENUKA:`nNo Response") | Out-Null
$ws.Range("I8").AddComment("RENUKA:`nNo Response") | Out-Null

# Update the active selection to the new last cell, as left by the editor
$ws.Range("I8").Select() | Out-Null
